# Generate Report for Handoff
# Updates the localization-status workbook: the 8cace8e5 file got a new
# handoff generated (fresh "Latest Handoff Datetime" / "Latest HO Xliff
# Generate Date" + an out-of-date-handback error), and both tracked files'
# Status moved from "Handed back: in sync with en-US" to "Ready for handoff".

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20d32e11806cd7185faffc79bea809b0db1d257e/e2e/8cace8e5-369e-4976-bf23-f08e7f682060.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b87bf1cf83f84e7fc407ea8c9742322747ec4585/e2e/8cace8e5-369e-4976-bf23-f08e7f682060.md."

# --- Overview sheet ---
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-22 18:51:34"

# --- zh-cn sheet ---
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-22 18:51:28"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-22 18:51:34"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
